$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 3).Value = 69579
$ws.Cells.Item(9, 5).Value = 191683877
$ws.Cells.Item(10, 3).Value = 278214
$ws.Cells.Item(10, 5).Value = 1752681420
$ws.Cells.Item(17, 3).Value = 134754
$ws.Cells.Item(17, 5).Value = 296820206
$ws.Cells.Item(38, 3).Value = 27073
$ws.Cells.Item(38, 5).Value = 83546230
$ws.Cells.Item(54, 3).Value = 17837
$ws.Cells.Item(54, 5).Value = 32098865
$ws.Cells.Item(69, 3).Value = 20737
$ws.Cells.Item(69, 5).Value = 62181544
$ws.Cells.Item(99, 3).Value = 136581
$ws.Cells.Item(99, 5).Value = 863240545
$ws.Cells.Item(126, 3).Value = 5647
$ws.Cells.Item(126, 5).Value = 8181370
$ws.Cells.Item(169, 3).Value = 562674
$ws.Cells.Item(169, 5).Value = 1286269320
$ws.Cells.Item(170, 3).Value = 367599
$ws.Cells.Item(170, 5).Value = 2848423734
$ws.Cells.Item(174, 3).Value = 357393
$ws.Cells.Item(174, 5).Value = 1020291034
$ws.Cells.Item(175, 3).Value = 125703
$ws.Cells.Item(175, 5).Value = 815960902
$ws.Cells.Item(179, 3).Value = 235816
$ws.Cells.Item(179, 5).Value = 813759583
$ws.Cells.Item(203, 3).Value = 13108
$ws.Cells.Item(203, 5).Value = 33026608
$ws.Cells.Item(205, 3).Value = 11135
$ws.Cells.Item(205, 5).Value = 44608881
$ws.Cells.Item(243, 3).Value = 28204
$ws.Cells.Item(243, 5).Value = 90602697
$ws.Cells.Item(257, 3).Value = 182554
$ws.Cells.Item(257, 5).Value = 1063849900
$ws.Cells.Item(262, 3).Value = 38988
$ws.Cells.Item(262, 5).Value = 124760356
$ws.Cells.Item(266, 3).Value = 71667
$ws.Cells.Item(266, 5).Value = 219458739
$ws.Cells.Item(311, 3).Value = 190864
$ws.Cells.Item(311, 5).Value = 586829007
$ws.Cells.Item(323, 3).Value = 94726
$ws.Cells.Item(323, 5).Value = 178876673
